$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two tasks that were in "To do" (D4) and "Doing" (E4) are now finished,
# so clear those in-progress cells...
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# ...and record them as completed in the "Done" column instead.
$ws.Range("G10").Value = "Compreender o código do projeto (José Pereira)"
$ws.Range("G11").Value = "Jogar FreeCol (José Pereira)"

# Leave the view the way the user had it when they finished: scrolled to
# show column B on the left, with the newly-added Done cells selected.
$ws.Activate()
$ws.Range("G10:G11").Select()
$excel.ActiveWindow.ScrollColumn = 2   # topLeftCell -> column B
$excel.ActiveWindow.ScrollRow = 1      # topLeftCell -> row 1
